$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 30
$ws.Cells.Item(30, 1).Value = 42603.657893518517
$ws.Cells.Item(30, 2).Value = "Noun"
$ws.Cells.Item(30, 3).Value = 3320
$ws.Cells.Item(30, 4).Value = 326
$ws.Cells.Item(30, 5).Value = 31
$ws.Cells.Item(30, 6).Value = 4
$ws.Cells.Item(30, 7).Value = 17
$ws.Cells.Item(30, 8).Value = 19
$ws.Cells.Item(30, 9).Value = 80
$ws.Cells.Item(30, 10).Value = 0
$ws.Cells.Item(30, 11).Value = 4
$ws.Cells.Item(30, 12).Value = 0
$ws.Cells.Item(30, 13).Value = 100

# Row 31
$ws.Cells.Item(31, 1).Value = 42603.690972222219
$ws.Cells.Item(31, 2).Value = "Noun"
$ws.Cells.Item(31, 3).Value = 3113
$ws.Cells.Item(31, 4).Value = 326
$ws.Cells.Item(31, 5).Value = 31
$ws.Cells.Item(31, 6).Value = 4
$ws.Cells.Item(31, 7).Value = 17
$ws.Cells.Item(31, 8).Value = 19
$ws.Cells.Item(31, 9).Value = 80
$ws.Cells.Item(31, 10).Value = 0
$ws.Cells.Item(31, 11).Value = 4
$ws.Cells.Item(31, 12).Value = 0
$ws.Cells.Item(31, 13).Value = 100

# Row 32
$ws.Cells.Item(32, 1).Value = 42603.691527777781
$ws.Cells.Item(32, 2).Value = "Noun"
$ws.Cells.Item(32, 3).Value = 3174
$ws.Cells.Item(32, 4).Value = 326
$ws.Cells.Item(32, 5).Value = 31
$ws.Cells.Item(32, 6).Value = 4
$ws.Cells.Item(32, 7).Value = 17
$ws.Cells.Item(32, 8).Value = 19
$ws.Cells.Item(32, 9).Value = 80
$ws.Cells.Item(32, 10).Value = 0
$ws.Cells.Item(32, 11).Value = 4
$ws.Cells.Item(32, 12).Value = 0
$ws.Cells.Item(32, 13).Value = 100

# Row 33
$ws.Cells.Item(33, 1).Value = 42603.692523148151
$ws.Cells.Item(33, 2).Value = "Noun"
$ws.Cells.Item(33, 3).Value = 3080
$ws.Cells.Item(33, 4).Value = 326
$ws.Cells.Item(33, 5).Value = 31
$ws.Cells.Item(33, 6).Value = 4
$ws.Cells.Item(33, 7).Value = 17
$ws.Cells.Item(33, 8).Value = 19
$ws.Cells.Item(33, 9).Value = 80
$ws.Cells.Item(33, 10).Value = 0
$ws.Cells.Item(33, 11).Value = 4
$ws.Cells.Item(33, 12).Value = 0
$ws.Cells.Item(33, 13).Value = 100

# Row 34
$ws.Cells.Item(34, 1).Value = 42603.692824074074
$ws.Cells.Item(34, 2).Value = "Noun"
$ws.Cells.Item(34, 3).Value = 3047
$ws.Cells.Item(34, 4).Value = 326
$ws.Cells.Item(34, 5).Value = 31
$ws.Cells.Item(34, 6).Value = 4
$ws.Cells.Item(34, 7).Value = 17
$ws.Cells.Item(34, 8).Value = 19
$ws.Cells.Item(34, 9).Value = 80
$ws.Cells.Item(34, 10).Value = 0
$ws.Cells.Item(34, 11).Value = 4
$ws.Cells.Item(34, 12).Value = 0
$ws.Cells.Item(34, 13).Value = 100

# Row 35
$ws.Cells.Item(35, 1).Value = 42603.693449074075
$ws.Cells.Item(35, 2).Value = "Noun"
$ws.Cells.Item(35, 3).Value = 3182
$ws.Cells.Item(35, 4).Value = 326
$ws.Cells.Item(35, 5).Value = 31
$ws.Cells.Item(35, 6).Value = 4
$ws.Cells.Item(35, 7).Value = 17
$ws.Cells.Item(35, 8).Value = 19
$ws.Cells.Item(35, 9).Value = 80
$ws.Cells.Item(35, 10).Value = 0
$ws.Cells.Item(35, 11).Value = 4
$ws.Cells.Item(35, 12).Value = 0
$ws.Cells.Item(35, 13).Value = 100
